$wb = $excel.ActiveWorkbook

# 1) Remove the "Wine vocabulary" sheet entirely
$wb.Worksheets.Item("Wine vocabulary").Delete() | Out-Null

# 2) Populate "Wine tasting" with the expanded flash-card deck
$ws = $wb.Worksheets.Item("Wine tasting")

$ws.Range("A5").Value = '1. Look'
$ws.Range("B5").Value = 'The clarity and the color'
$ws.Range("A6").Value = 'Clarity'
$ws.Range("B6").Value = 'Brilliant - Bright - Clear - Dull - Hazy - Cloudy / Sediment? Crystals?'
$ws.Range("A7").Value = 'Color (White wine)'
$ws.Range("B7").Value = 'dilute, v. pale straw, pale straw, medium straw, yellow, lt. gold, medium gold, amber, tawny, adobe, brown'
$ws.Range("A8").Value = 'Color (Red wine)'
$ws.Range("B8").Value = 'light / medium / dark: garnet, ruby red, purplish-red, or combinations / tawny, brown edge, any rim variation?'
$ws.Range("A9").Value = 'Old color'
$ws.Range("B9").Value = 'The older is a wine (white or red), the closer to the brown will the color be'
$ws.Range("A10").Value = '2. Smell'
$ws.Range("B10").Value = 'First sniff, swirling and Second sniff'
$ws.Range("A11").Value = 'First sniff'
$ws.Range("B11").Value = 'Is the wine presentable? Fresh nose or aged-developed bouquet? OK to proceed? Any FAULTS? Corked, vinegary, oxidized / stale, stinky / eggy, sulfites'
$ws.Range("A12").Value = 'Swirling'
$ws.Range("B12").Value = 'Needs a good wine glass, makes the flavor easier to smell'
$ws.Range("A13").Value = 'Second sniff'
$ws.Range("B13").Value = 'Fruity / floral / spicy / berry / ripeness-level: under ripe through jammy, raisiny / estery / herbaceous / vegetal / grassy / nutty / minerally? Also think about intensity. Woody: fresh woody, vanilla, caramel, brown sugar, smoky-char,coconut '
$ws.Range("A14").Value = '3. Taste'
$ws.Range("B14").Value = 'Find the basic tastes, the retro odors, and identify the mouthfeel (body), and mesure the persistence'
$ws.Range("A15").Value = 'Basic Tastes'
$ws.Range("B15").Value = 'Sweetness, Acidity (sourness), Bitterness (think about balance)'
$ws.Range("A16").Value = 'Retro-Odors'
$ws.Range("B16").Value = 'Do the retro-nasal odors while on your palate confirm the ortho-nasal odors? Any new odors?'
$ws.Range("A17").Value = 'Mouthfeel / Body'
$ws.Range("B17").Value = 'The wine''s weight (lt. body, medium body, full body) Mouthfeel: smooth or rough, astringent? Hotness from alcohol?'
$ws.Range("A18").Value = 'Persistence'
$ws.Range("B18").Value = 'Short (under 30 sec), medium (up to a minute), or long finish (a few minutes)'
$ws.Range("A19").Value = '4. Feeling'
$ws.Range("B19").Value = 'Did you like it? How much would you rank it (compare to other similar ones)'
$ws.Range("A20").Value = 'Two basic rules'
$ws.Range("B20").Value = 'Take your time. Be attentive'
$ws.Range("A21").Value = 'Blind tasting'
$ws.Range("B21").Value = 'Tasting without seeing the bottle, and without knowing the price'
$ws.Range("A22").Value = 'Tears'
$ws.Range("B22").Value = 'The tears you see from the wine does not mean that the wine is good or bad'
$ws.Range("A23").Value = 'Tricks for smelling'
$ws.Range("B23").Value = 'Don''t wear a perfume, train identifying every flavor to train your nose'
$ws.Range("A24").Value = 'Tricks for looking'
$ws.Range("B24").Value = 'Put a white background, a paper is good, and a reflection of sun light would be the best'
$ws.Range("A25").Value = 'Sweet'
$ws.Range("B25").Value = 'It is the opposite of dry. If you eat a grape and you remove the skin you will find what sweet is. Be careful sweet doesn''t mean fruity.'
$ws.Range("A26").Value = 'Acidity'
$ws.Range("B26").Value = 'The sourness of the wine. In white wine it compensate the sweet taste. It generates saliva in your mouth.'
$ws.Range("A27").Value = 'Tanin'
$ws.Range("B27").Value = 'The bitterness of the red wine, it comes from the skin. If you eat the skin only of a grape you will find this taste. It generate the dry taste.'
$ws.Range("A28").Value = 'Balance'
$ws.Range("B28").Value = 'A wine is balanced if sweet, acidity, tanin, and alcohol are compensating each other'
$ws.Range("A29").Value = 'Lenght'
$ws.Range("B29").Value = 'The way the wine is doing in your palate'
$ws.Range("A30").Value = 'Depth'
$ws.Range("B30").Value = 'The layers of taste and aromas of the wine'
$ws.Range("A31").Value = 'Complexity'
$ws.Range("B31").Value = 'A wine is complex if it has a good lenght and depth'
$ws.Range("A32").Value = 'Typicality'
$ws.Range("B32").Value = 'If the tastes are similar to the other ones of the same race'
$ws.Range("A33").Value = 'What is a bad wine'
$ws.Range("B33").Value = 'If it taste like a rotten fruit, like vinegar, chemical, oxyde, burnt or cork'

# 3) Restore view state: scrolled so row 9 is at the top, active cell B33 selected
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 9 | Out-Null
$ws.Range("B33").Select() | Out-Null

